$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old "faciliteiten" text (column E) to new abbreviated codes
$map = @{
    "gymastiek" = "LO"
    "latijn|grieks|algemeen" = "LA|GR|ALG"
    "muziek" = "MU"
    "algemeen" = "ALG"
    "computer" = "CP"
    "geschiedenis|algemeen" = "GS|ALG"
    "wiskunde|algemeen" = "WI|ALG"
    "aardrijkskunde|algemeen" = "AK|ALG"
    "biologie" = "BI"
    "scheikunde|natuurkunde" = "NA|SK"
}

for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $old = $cell.Value()
    if ($map.ContainsKey($old)) {
        $cell.Value = $map[$old]
    }
}

# Update the active selection shown in the worksheet view
$ws.Range("B17").Select()
